# Appends the new "ਕਾਹੇ ਕੀ ਕੁਸਲਾਤ ਹਾਥਿ ਦੀਪੁ ਕੂਏ ਪਰੈ ॥੨੧੬॥" word-by-word assessment
# as rows 241-247 (verse ends with a new Shalok-End row), per the new verse
# assessment commit. Existing rows 1-240 are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 241
$ws.Range("A241").Value = 'ਕਾਹੇ ਕੀ ਕੁਸਲਾਤ ਹਾਥਿ ਦੀਪੁ ਕੂਏ ਪਰੈ ॥੨੧੬॥'
$ws.Range("B241").Value = 'Why (speak) of well-being, (when) the lamp in hand falls into the well?'
$ws.Range("D241").Value = 'ਕਾਹੇ'
$ws.Range("E241").Value = 'ਕਿਸ ਵਾਸਤੇ?'
$ws.Range("F241").Value = 'ਕਾਹੇ'
$ws.Range("H241").Value = 'Reason / ਕਾਰਣ ਵਾਚਕ'
$ws.Range("K241").Value = 'Adverb / ਕਿਰਿਆ ਵਿਸੇਸ਼ਣ'
$ws.Range("L241").Value = 1
$ws.Range("M241").Value = 0
$ws.Range("N241").Value = 54731
$ws.Range("O241").Value = 2
$ws.Range("P241").Value = 1
$ws.Range("Q241").Value = 3818
$ws.Range("R241").Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Range("T241").Value = 'ਕਬੀਰ ਜੀ'
$ws.Range("Y241").Value = 'ਸ਼ਲੋਕ'
$ws.Range("Z241").Value = 'Shalok End'
$ws.Range("AA241").Value = 1376
$ws.Range("AB241").Value = 0
$ws.Range("AC241").Value = 0

# Row 242
$ws.Range("A242").Value = 'ਕਾਹੇ ਕੀ ਕੁਸਲਾਤ ਹਾਥਿ ਦੀਪੁ ਕੂਏ ਪਰੈ ॥੨੧੬॥'
$ws.Range("B242").Value = 'Why (speak) of well-being, (when) the lamp in hand falls into the well?'
$ws.Range("D242").Value = 'ਕੀ'
$ws.Range("E242").Value = 'ਦੀ'
$ws.Range("F242").Value = 'ਕੀ'
$ws.Range("G242").Value = 'Singular / ਇਕ'
$ws.Range("H242").Value = 'Of'
$ws.Range("I242").Value = 'Feminine / ਇਸਤਰੀ'
$ws.Range("K242").Value = 'Postposition / ਸੰਬੰਧਕ'
$ws.Range("L242").Value = 1
$ws.Range("M242").Value = 1
$ws.Range("N242").Value = 54731
$ws.Range("O242").Value = 2
$ws.Range("P242").Value = 1
$ws.Range("Q242").Value = 3818
$ws.Range("R242").Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Range("T242").Value = 'ਕਬੀਰ ਜੀ'
$ws.Range("Y242").Value = 'ਸ਼ਲੋਕ'
$ws.Range("Z242").Value = 'Shalok End'
$ws.Range("AA242").Value = 1376
$ws.Range("AB242").Value = 0
$ws.Range("AC242").Value = 0

# Row 243
$ws.Range("A243").Value = 'ਕਾਹੇ ਕੀ ਕੁਸਲਾਤ ਹਾਥਿ ਦੀਪੁ ਕੂਏ ਪਰੈ ॥੨੧੬॥'
$ws.Range("B243").Value = 'Why (speak) of well-being, (when) the lamp in hand falls into the well?'
$ws.Range("D243").Value = 'ਕੁਸਲਾਤ'
$ws.Range("E243").Value = 'ਸੁਖ'
$ws.Range("F243").Value = 'ਮੁਕਤਾ'
$ws.Range("G243").Value = 'Singular / ਇਕ'
$ws.Range("H243").Value = 'Genitive ਸੰਬੰਧ'
$ws.Range("I243").Value = 'Feminine / ਇਸਤਰੀ'
$ws.Range("J243").Value = 'ਮੁਕਤਾ Ending'
$ws.Range("K243").Value = 'Noun / ਨਾਂਵ'
$ws.Range("L243").Value = 1
$ws.Range("M243").Value = 2
$ws.Range("N243").Value = 54731
$ws.Range("O243").Value = 2
$ws.Range("P243").Value = 1
$ws.Range("Q243").Value = 3818
$ws.Range("R243").Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Range("T243").Value = 'ਕਬੀਰ ਜੀ'
$ws.Range("Y243").Value = 'ਸ਼ਲੋਕ'
$ws.Range("Z243").Value = 'Shalok End'
$ws.Range("AA243").Value = 1376
$ws.Range("AB243").Value = 0
$ws.Range("AC243").Value = 0

# Row 244
$ws.Range("A244").Value = 'ਕਾਹੇ ਕੀ ਕੁਸਲਾਤ ਹਾਥਿ ਦੀਪੁ ਕੂਏ ਪਰੈ ॥੨੧੬॥'
$ws.Range("B244").Value = 'Why (speak) of well-being, (when) the lamp in hand falls into the well?'
$ws.Range("D244").Value = 'ਹਾਥਿ'
$ws.Range("E244").Value = 'ਹੱਥ ਵਿਚ'
$ws.Range("F244").Value = 'ਿ'
$ws.Range("G244").Value = 'Singular / ਇਕ'
$ws.Range("H244").Value = 'Locative ਅਧਿਕਰਣ'
$ws.Range("I244").Value = 'Masculine / ਪੁਲਿੰਗ'
$ws.Range("J244").Value = 'ਮੁਕਤਾ Ending'
$ws.Range("K244").Value = 'Noun / ਨਾਂਵ'
$ws.Range("L244").Value = 1
$ws.Range("M244").Value = 3
$ws.Range("N244").Value = 54731
$ws.Range("O244").Value = 2
$ws.Range("P244").Value = 1
$ws.Range("Q244").Value = 3818
$ws.Range("R244").Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Range("T244").Value = 'ਕਬੀਰ ਜੀ'
$ws.Range("Y244").Value = 'ਸ਼ਲੋਕ'
$ws.Range("Z244").Value = 'Shalok End'
$ws.Range("AA244").Value = 1376
$ws.Range("AB244").Value = 0
$ws.Range("AC244").Value = 0

# Row 245
$ws.Range("A245").Value = 'ਕਾਹੇ ਕੀ ਕੁਸਲਾਤ ਹਾਥਿ ਦੀਪੁ ਕੂਏ ਪਰੈ ॥੨੧੬॥'
$ws.Range("B245").Value = 'Why (speak) of well-being, (when) the lamp in hand falls into the well?'
$ws.Range("D245").Value = 'ਦੀਪੁ'
$ws.Range("E245").Value = 'ਹਾਥਿ ਦੀਪੁ: ਹੱਥਾਂ ਵਿਚ ਦੀਵਾ'
$ws.Range("F245").Value = 'ੁ'
$ws.Range("G245").Value = 'Singular / ਇਕ'
$ws.Range("H245").Value = 'Nominative ਕਰਤਾ'
$ws.Range("I245").Value = 'Masculine / ਪੁਲਿੰਗ'
$ws.Range("J245").Value = 'ਮੁਕਤਾ Ending'
$ws.Range("K245").Value = 'Noun / ਨਾਂਵ'
$ws.Range("L245").Value = 1
$ws.Range("M245").Value = 4
$ws.Range("N245").Value = 54731
$ws.Range("O245").Value = 2
$ws.Range("P245").Value = 1
$ws.Range("Q245").Value = 3818
$ws.Range("R245").Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Range("T245").Value = 'ਕਬੀਰ ਜੀ'
$ws.Range("Y245").Value = 'ਸ਼ਲੋਕ'
$ws.Range("Z245").Value = 'Shalok End'
$ws.Range("AA245").Value = 1376
$ws.Range("AB245").Value = 0
$ws.Range("AC245").Value = 0

# Row 246
$ws.Range("A246").Value = 'ਕਾਹੇ ਕੀ ਕੁਸਲਾਤ ਹਾਥਿ ਦੀਪੁ ਕੂਏ ਪਰੈ ॥੨੧੬॥'
$ws.Range("B246").Value = 'Why (speak) of well-being, (when) the lamp in hand falls into the well?'
$ws.Range("D246").Value = 'ਕੂਏ'
$ws.Range("E246").Value = 'ਖੂਹ ਵਿਚ'
$ws.Range("F246").Value = 'ਕੂਏ'
$ws.Range("G246").Value = 'Singular / ਇਕ'
$ws.Range("H246").Value = 'Locative ਅਧਿਕਰਣ'
$ws.Range("I246").Value = 'Masculine / ਪੁਲਿੰਗ'
$ws.Range("J246").Value = 'ਕੰਨਾ Ending'
$ws.Range("K246").Value = 'Noun / ਨਾਂਵ'
$ws.Range("L246").Value = 1
$ws.Range("M246").Value = 5
$ws.Range("N246").Value = 54731
$ws.Range("O246").Value = 2
$ws.Range("P246").Value = 1
$ws.Range("Q246").Value = 3818
$ws.Range("R246").Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Range("T246").Value = 'ਕਬੀਰ ਜੀ'
$ws.Range("Y246").Value = 'ਸ਼ਲੋਕ'
$ws.Range("Z246").Value = 'Shalok End'
$ws.Range("AA246").Value = 1376
$ws.Range("AB246").Value = 0
$ws.Range("AC246").Value = 0

# Row 247
$ws.Range("A247").Value = 'ਕਾਹੇ ਕੀ ਕੁਸਲਾਤ ਹਾਥਿ ਦੀਪੁ ਕੂਏ ਪਰੈ ॥੨੧੬॥'
$ws.Range("B247").Value = 'Why (speak) of well-being, (when) the lamp in hand falls into the well?'
$ws.Range("D247").Value = 'ਪਰੈ'
$ws.Range("E247").Value = 'ਪੜੈ| ਪੈਂਦਾ ਹੈ'
$ws.Range("F247").Value = 'ਕਰੈ'
$ws.Range("G247").Value = 'Singular / ਇਕ'
$ws.Range("H247").Value = 'Present ਵਰਤਮਾਨ'
$ws.Range("I247").Value = 'Trans / ਨਪੁਂਸਕ'
$ws.Range("J247").Value = '3rd Person / ਅਨਯ ਪੁਰਖ'
$ws.Range("K247").Value = 'Verb / ਕਿਰਿਆ'
$ws.Range("L247").Value = 1
$ws.Range("M247").Value = 6
$ws.Range("N247").Value = 54731
$ws.Range("O247").Value = 2
$ws.Range("P247").Value = 1
$ws.Range("Q247").Value = 3818
$ws.Range("R247").Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Range("T247").Value = 'ਕਬੀਰ ਜੀ'
$ws.Range("Y247").Value = 'ਸ਼ਲੋਕ'
$ws.Range("Z247").Value = 'Shalok End'
$ws.Range("AA247").Value = 1376
$ws.Range("AB247").Value = 0
$ws.Range("AC247").Value = 0
